$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update survey response values
$ws.Range("D1").Value = 13
$ws.Range("B2").Value = 7
$ws.Range("A3").Value = 5
$ws.Range("A4").Value = 10
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 4

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("F10").Select()
